$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp banner (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 21 de Mayo de 2020 a las 11:35"

# --- Row 20: Belgica - refreshed case counts ---
$ws.Range("B20").Value = 56235
$ws.Range("C20").Value = 252
$ws.Range("D20").Value = 14988
$ws.Range("E20").Value = 32061
$ws.Range("G20").Value = 36
$ws.Range("H20").Value = 9186

# --- Row 63: Finlandia - refreshed case counts ---
$ws.Range("B63").Value = 6493
$ws.Range("C63").Value = 50
$ws.Range("E63").Value = 1389

# --- Rows 197/198: Nueva Caledonia <-> Santa Lucia swap position ---
$ws.Range("A197").Value = "Santa Lucia"
$ws.Range("A198").Value = "Nueva Caledonia"

# --- Rows 209/210: Seychelles <-> Montserrat swap (with their Recuperados/Muertes) ---
$ws.Range("A209").Value = "Montserrat"
$ws.Range("D209").Value = 10
$ws.Range("H209").Value = 1

$ws.Range("A210").Value = "Seychelles"
$ws.Range("D210").Value = 11
$ws.Range("H210").Value = 0

# --- Rows 214/215: Sahara Occidental <-> Bonaire, San Eustaquio y Saba swap position ---
$ws.Range("A214").Value = "Bonaire, San Eustaquio y Saba"
$ws.Range("A215").Value = "Sahara Occidental"
